# Lesson 12 wordlist reorder:
# rows 2-40 (verbs/phrases) moved after row 34 (new block of adjectives/body terms),
# and the remaining rows re-sequenced accordingly (per upstream word-list reshuffle).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (English) and Column B (Japanese) in final row order, row 1 = header.
$english = @(
    'English',
    'leg; foot',
    'meaning',
    'restroom',
    'stomach',
    'cold',
    'girlfriend',
    'boyfriend',
    'temperature (weather)',
    'cloudy weather',
    'match; game',
    'juice',
    'politics',
    'grade (on a test, etc.)',
    'cough',
    'throat',
    'tooth',
    'flower',
    'sunny weather',
    'clothes',
    'hangover',
    'present',
    'homesickness',
    'minus',
    'thing (concrete object)',
    'snow',
    'business to take care of',
    'sweet',
    'hurt; painful',
    'there are many',
    'narrow; not spacious',
    'inconvenient; to have a scheduling conflict',
    'bad',
    'nice',
    'to walk',
    'to catch a cold',
    'to be interested (in...)',
    'to lose',
    'to have a fever',
    'to become thirsty',
    'to cough',
    'to break up; to separate',
    'to get nervous',
    'to worry',
    'Get well soon.',
    'don''t look well',
    'probably; maybe',
    'as much as possible',
    'probably; ..., right?',
    'degrees (temperature)',
    'for two to three days',
    'because...',
    'for the first time',
    'very soon; in a few moments/days',
    'I have diarrhea',
    'I am constipated.',
    'I have my period.',
    'I have hay fever.',
    'I have an allergy to...',
    'I have a bad tooth.',
    'I sneeze.',
    'I have a runny nose.',
    'My back itches.',
    'I have rashes.',
    'I feel dizzy.',
    'I threw up.',
    'I am not feeling well.',
    'I burned myself.',
    'I broke my leg.',
    'I hurt myself.',
    'physician',
    'dermatologist',
    'surgeon',
    'obstetrician and gynecologist',
    'orthopedic surgeon',
    'ophthalmologist',
    'dentist',
    'otorhinolaryngologist; ENT doctor',
    'antibiotic',
    'X-ray',
    'operation',
    'injection',
    'thermometer',
    'old times',
    'old tale',
    'ancient times',
    'once upon a time',
    'people',
    'sometimes',
    'various',
    'God',
    'shrine',
    'Shinto religion',
    'Kobe City',
    'early',
    'to get up early',
    'early morning',
    'to get up',
    'to wake someone up',
    'to stand up',
    'cow',
    'milk',
    'beef',
    'calf; veal',
    'to use',
    'ambassador',
    '"Occupied"',
    'errand',
    'to work',
    'both husband and wife working',
    'labor',
    'to bring (a person) back',
    'United Nations',
    'consecutive holidays',
    'to separate',
    'not in particular',
    'special',
    'discrimination',
    'separately',
    'once',
    'near future',
    'temperature',
    '30 degrees',
    'attitude',
    'red color',
    'red',
    'baby',
    'the equator',
    'the Red Cross',
    'blue color',
    'blue',
    'youth',
    'blue sky',
    'green light',
    'color',
    'various',
    'scenery',
    'characteristic'
)

$japanese = @(
    'Japanese',
    '足|あし',
    '意味|いみ',
    'お手洗い|おてあらい',
    'おなか',
    '風邪|かぜ',
    '彼女|かのじょ',
    '彼|かれ',
    '気温|きおん',
    '曇り|くもり',
    '試合|しあい',
    'ジュース',
    '政治|せいじ',
    '成績|せいせき',
    'せき',
    'のど',
    '歯|は',
    '花|はな',
    '晴れ|はれ',
    '服|ふく',
    '二日酔い|ふつかよい',
    'プレゼント',
    'ホームシック',
    'マイナス',
    '物|もの',
    '雪|ゆき',
    '用事|ようじ',
    '甘い|あまい',
    '痛い|いたい',
    '多い|おおい',
    '狭い|せまい',
    '都合が悪い|つごうがわるい',
    '悪い|わるい',
    '素敵|すてき（な）',
    '歩く|あるく',
    '風邪をひく|かぜをひく',
    '興味がある|きょうみがある',
    'なくす',
    '熱がある|ねつがある',
    'のどが渇く|のどがかわく',
    'せきが出る|せきがでる',
    '別れる|わかれる',
    '緊張する|きんちょうする',
    '心配する|しんぱいする',
    'お大事に|おだいじに',
    '元気がない|げんきがない',
    '多分|たぶん',
    'できるだけ',
    '～でしょう',
    '～度|～ど',
    '二三日|にさんにち',
    '～ので',
    '初めて|はじめて',
    'もうすぐ',
    '下痢です。|げりです。',
    '便秘です。|べんぴです。',
    '生理です。|せいりです。',
    '花粉症です。|かふんしょうです。',
    '（～に）アレルギーがあります。',
    '虫歯があります。|むしばがあります。',
    'くしゃみが出ます。|くしゃみがでます。',
    '鼻水が出ます。|はなみずがでます。',
    '背中がかゆいです。|せなかがかゆいです。',
    '発疹があります。|はっしんがあります。',
    'めまいがします。',
    '吐きました。|はきました。',
    '気分が悪いです。|きぶんがわるいです。',
    'やけどをしました。',
    '足の骨を折りました。|あしのほねをおりました。',
    'けがをしました。',
    '内科|ないか',
    '皮膚科|ひふか',
    '外科|げか',
    '産婦人科|さんふじんか',
    '整形外科|せいけいげか',
    '眼科|がんか',
    '歯科|しか',
    '耳鼻科|じびか',
    '抗生物質|こうせいぶっしつ',
    'レントゲン',
    '手術|しゅじゅつ',
    '注射|ちゅうしゃ',
    '体温計|たいおんけい',
    '昔|むかし',
    '昔話|むかしばなし',
    '大昔|おおむかし',
    '昔々|むかしむかし',
    '人々|ひとびと',
    '時々|ときどき',
    '色々な|いろいろな',
    '神さま|かみさま',
    '神社|じんじゃ',
    '神道|しんとう',
    '神戸市|こうべし',
    '早い|はやい',
    '早起きする|はやおきする',
    '早朝|そうちょう',
    '起きる|おきる',
    '起こす|おこす',
    '起立する|きりつする',
    '牛|うし',
    '牛乳|ぎゅうにゅう',
    '牛肉|ぎゅうにく',
    '子牛|こうし',
    '使う|つかう',
    '大使|たいし',
    '使用中|しようちゅう',
    'お使い|おつかい',
    '働く|はたらく',
    '共働き|ともばたらき',
    '労働|ろうどう',
    '連れて帰る|つれてかえる',
    '国連|こくれん',
    '連休|れんきゅう',
    '別れる|わかれる',
    '別に|べつに',
    '特別な|とくべつな',
    '差別|さべつ',
    '別々に|べつべつに',
    '一度|いちど',
    '今度|こんど',
    '温度|おんど',
    '三十度|さんじゅうど',
    '態度|たいど',
    '赤|あか',
    '赤い|あかい',
    '赤ちゃん|あかちゃん',
    '赤道|せきどう',
    '赤十字|せきじゅうじ',
    '青|あお',
    '青い|あおい',
    '青年|せいねん',
    '青空|あおぞら',
    '青信号|あおしんごう',
    '色|いろ',
    '色々な|いろいろな',
    '景色|けしき',
    '特色|とくしょく'
)

for ($i = 0; $i -lt $english.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $english[$i]
    $ws.Cells.Item($row, 2).Value = $japanese[$i]
}
